# res mob sim mod
# Adds the "percent change" rows (49-61) to Sheet1: six rows of
# Scenario 1 figures, a blank separator row, then six rows of Scenario 2
# figures. Styling is picked up from existing same-purpose cells via
# copy/paste-format so the new rows match the workbook's established look
# (label column uses the B2-style "text" format, most numeric columns use
# the D2-style "0" format, and F49/G49 reuse the bold/centered C42-style
# "0" format used for the other scenario-total header rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- seed formats for the new rows by copying from representative cells ---

# Label cells (column B) -> same look as existing row labels (e.g. B2)
# (row 55 is a blank separator row with no label, so it is skipped here)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B49:B54").PasteSpecial(-4122) | Out-Null
$ws.Range("B56:B61").PasteSpecial(-4122) | Out-Null

# Generic numeric cells -> same look as D2 ("0" format, Arial, centered)
$ws.Range("D2").Copy() | Out-Null
$ws.Range("C49:G61").PasteSpecial(-4122) | Out-Null

# F49/G49 use the bolder "0" format shared with the other scenario totals
$ws.Range("C42").Copy() | Out-Null
$ws.Range("F49:G49").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- row 49: Scenario 1 totals -------------------------------------------
$ws.Range("B49").Value = "scenario1_totaldwellings_percchange"
$ws.Range("C49").Value = 403.517587939699
$ws.Range("D49").Value = 372.080979284369
$ws.Range("E49").Value = 112.166003411029
$ws.Range("F49").Value = 48.5584444609998
$ws.Range("G49").Value = 1448.35309617918

# --- row 50: Scenario 1 avg bedrooms --------------------------------------
$ws.Range("B50").Value = "scenario1_avgbedrooms_percchange"
$ws.Range("C50").Value = -35.7142857142857
$ws.Range("D50").Value = -32.5287356321839
$ws.Range("E50").Value = -19.2374956278419
$ws.Range("F50").Value = -0.74224021592442
$ws.Range("G50").Value = -33.8188347964884

# --- row 51: Scenario 1 single detached -----------------------------------
$ws.Range("B51").Value = "scenario1_singledetached_percchange"
$ws.Range("C51").Value = -12.6578840146681
$ws.Range("D51").Value = -2.48804233307282
$ws.Range("E51").Value = 12.1708844083405
$ws.Range("F51").Value = -4.05894169177998
$ws.Range("G51").Value = 0.149413687478176

# --- row 52: Scenario 1 row/town/semi --------------------------------------
$ws.Range("B52").Value = "scenario1_rowtownsemi_percchange"
$ws.Range("C52").Value = -12.6578840146679
$ws.Range("D52").Value = 319.649112661446
$ws.Range("E52").Value = 535.535104319977
$ws.Range("F52").Value = -4.0589416917797
$ws.Range("G52").Value = 0.149413687478025

# --- row 53: Scenario 1 apt small ------------------------------------------
$ws.Range("B53").Value = "scenario1_aptsmall_percchange"
$ws.Range("C53").Value = 306.65833487615
$ws.Range("D53").Value = 67.7296930269655
$ws.Range("E53").Value = 30.2664073655311
$ws.Range("F53").Value = 6.6423992453538
$ws.Range("G53").Value = 0.149413687478062

# --- row 54: Scenario 1 apt large ------------------------------------------
$ws.Range("B54").Value = "scenario1_aptlarge_percchange"
$ws.Range("C54").Value = 561.461493204446
$ws.Range("D54").Value = 2029.68510736727
$ws.Range("E54").Value = -14.0672834060686
$ws.Range("F54").Value = 295.048906066792
$ws.Range("G54").Value = 17330.1032242954

# --- row 55: blank separator row (keeps the "0"-format styling only) ------
$ws.Range("C55").Value = ""
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = ""
$ws.Range("G55").Value = ""

# --- row 56: Scenario 2 totals ----------------------------------------------
$ws.Range("B56").Value = "scenario2_totaldwellings_percchange"
$ws.Range("C56").Value = 320.302594051338
$ws.Range("D56").Value = 355.461393596987
$ws.Range("E56").Value = 120.238772029562
$ws.Range("F56").Value = 47.6868062281159
$ws.Range("G56").Value = 1383.13570487484

# --- row 57: Scenario 2 avg bedrooms ----------------------------------------
$ws.Range("B57").Value = "scenario2_avgbedrooms_percchange"
$ws.Range("C57").Value = -15.9770037850021
$ws.Range("D57").Value = -32.3371647509579
$ws.Range("E57").Value = -21.860790486184
$ws.Range("F57").Value = 9.87854251012146
$ws.Range("G57").Value = -33.7190742218675

# --- row 58: Scenario 2 single detached --------------------------------------
$ws.Range("B58").Value = "scenario2_singledetached_percchange"
$ws.Range("C58").Value = -12.6582226922109
$ws.Range("D58").Value = -3.21714922595285
$ws.Range("E58").Value = 20.3036107618776
$ws.Range("F58").Value = -4.05897826018926
$ws.Range("G58").Value = 0.155994050642836

# --- row 59: Scenario 2 row/town/semi ----------------------------------------
$ws.Range("B59").Value = "scenario2_rowtownsemi_percchange"
$ws.Range("C59").Value = -12.6582226922108
$ws.Range("D59").Value = 305.255328402163
$ws.Range("E59").Value = 440.443697359411
$ws.Range("F59").Value = 246.886768607318
$ws.Range("G59").Value = 0.155994050643039

# --- row 60: Scenario 2 apt small ---------------------------------------------
$ws.Range("B60").Value = "scenario2_aptsmall_percchange"
$ws.Range("C60").Value = 455.700982445775
$ws.Range("D60").Value = 64.0220327936284
$ws.Range("E60").Value = 153.685124400053
$ws.Range("F60").Value = 5.29338704828814
$ws.Range("G60").Value = 1502.02277773789

# --- row 61: Scenario 2 apt large ---------------------------------------------
$ws.Range("B61").Value = "scenario2_aptlarge_percchange"
$ws.Range("C61").Value = 422.233363603296
$ws.Range("D61").Value = 1942.75361738141
$ws.Range("E61").Value = 4.88601968203052
$ws.Range("F61").Value = 140.94077941064
$ws.Range("G61").Value = 9599.56008349296

# --- refresh the sheet's scroll/selection position ------------------------
$ws.Activate() | Out-Null
$ws.Range("G52").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
